$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates scraped on Sun Feb 12 07:20:36 UTC 2023 (coinranking snapshot refresh).
# Columns D (Price) and E (Volume 1h) are numeric-looking text values in the source
# sheet (t="inlineStr"), so each is written with an explicit Text number format and
# then restored to the default "Normal" style to avoid Excel auto-coercing them into
# real numbers/percentages (and to keep the cell style index unchanged, s=0).
$priceAndVolumeUpdates = @(
    @{ Cell = 'D2'; Value = '307.38' }
    @{ Cell = 'E2'; Value = '-0.35%' }
    @{ Cell = 'D3'; Value = '41.03' }
    @{ Cell = 'E3'; Value = '0.52%' }
    @{ Cell = 'D4'; Value = '5.237' }
    @{ Cell = 'E4'; Value = '2.36%' }
    @{ Cell = 'D5'; Value = '0.07666' }
    @{ Cell = 'E5'; Value = '0.73%' }
    @{ Cell = 'D6'; Value = '1.633' }
    @{ Cell = 'E6'; Value = '0.49%' }
    @{ Cell = 'D7'; Value = '0.9158' }
    @{ Cell = 'E7'; Value = '1.55%' }
    @{ Cell = 'D9'; Value = '0.1246' }
    @{ Cell = 'E9'; Value = '13.73%' }
    @{ Cell = 'E10'; Value = '2.94%' }
    @{ Cell = 'D11'; Value = '0.09093' }
    @{ Cell = 'E11'; Value = '-1.06%' }
    @{ Cell = 'D12'; Value = '0.04158' }
    @{ Cell = 'E12'; Value = '-0.43%' }
    @{ Cell = 'E13'; Value = '-0.02%' }
    @{ Cell = 'D14'; Value = '0.001259' }
    @{ Cell = 'E14'; Value = '0.03%' }
    @{ Cell = 'D15'; Value = '0.005850' }
    @{ Cell = 'E15'; Value = '-0.21%' }
    @{ Cell = 'D17'; Value = '3.345' }
    @{ Cell = 'E17'; Value = '-0.22%' }
    @{ Cell = 'D18'; Value = '4.311' }
    @{ Cell = 'E18'; Value = '1.35%' }
    @{ Cell = 'D19'; Value = '0.3335' }
    @{ Cell = 'E19'; Value = '1.13%' }
    @{ Cell = 'D20'; Value = '7.313' }
    @{ Cell = 'E20'; Value = '11.55%' }
    @{ Cell = 'D21'; Value = '0.1384' }
    @{ Cell = 'E21'; Value = '1.72%' }
    @{ Cell = 'D22'; Value = '0.2715' }
    @{ Cell = 'E22'; Value = '1.25%' }
    @{ Cell = 'D23'; Value = '0.04083' }
    @{ Cell = 'E23'; Value = '0.41%' }
    @{ Cell = 'D24'; Value = '0.001264' }
    @{ Cell = 'E24'; Value = '3.32%' }
    @{ Cell = 'D25'; Value = '0.004281' }
    @{ Cell = 'E25'; Value = '4.63%' }
    @{ Cell = 'E26'; Value = '-2.15%' }
    @{ Cell = 'D38'; Value = '0.02489' }
    @{ Cell = 'E38'; Value = '4.84%' }
    @{ Cell = 'D39'; Value = '0.05315' }
    @{ Cell = 'E39'; Value = '2.47%' }
    @{ Cell = 'D40'; Value = '0.007846' }
    @{ Cell = 'E40'; Value = '0.67%' }
    @{ Cell = 'D41'; Value = '0.1314' }
    @{ Cell = 'E41'; Value = '0.96%' }
    @{ Cell = 'D42'; Value = '0.006893' }
    @{ Cell = 'E42'; Value = '2.11%' }
    @{ Cell = 'D44'; Value = '0.007655' }
    @{ Cell = 'E44'; Value = '-8.29%' }
    @{ Cell = 'D45'; Value = '0.3055' }
    @{ Cell = 'E45'; Value = '-0.71%' }
    @{ Cell = 'D46'; Value = '0.00006720' }
    @{ Cell = 'E46'; Value = '-3.66%' }
    @{ Cell = 'E47'; Value = '0.16%' }
    @{ Cell = 'D48'; Value = '0.4364' }
    @{ Cell = 'E48'; Value = '1,271.95%' }
    @{ Cell = 'D49'; Value = '0.003107' }
    @{ Cell = 'E50'; Value = '0.16%' }
    @{ Cell = 'D51'; Value = '0.0002005' }
    @{ Cell = 'E51'; Value = '0.16%' }
)

foreach ($update in $priceAndVolumeUpdates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
    $cell.Style = "Normal"
}

# Columns B (Coin) and C (Link) are plain (non-numeric-looking) text, so a direct
# assignment is safe and keeps the default style untouched.
$coinAndLinkUpdates = @(
    @{ Cell = 'B17'; Value = 'LEO' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'B18'; Value = 'GateToken' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'B19'; Value = 'BitpandaEcosystemToken' }
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best' }
    @{ Cell = 'B20'; Value = 'MCDex' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb' }
    @{ Cell = 'B21'; Value = 'ProBitToken' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob' }
    @{ Cell = 'B22'; Value = 'ZBToken' }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb' }
    @{ Cell = 'B23'; Value = 'CoinExToken' }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' }
    @{ Cell = 'B24'; Value = 'BitKan' }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan' }
    @{ Cell = 'B25'; Value = 'HotbitToken' }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb' }
)

foreach ($update in $coinAndLinkUpdates) {
    $ws.Range($update.Cell).Value = $update.Value
}
